$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at the very top for "Date and Time"
$ws.Rows("1:1").Insert()
$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-11 17:42:37.788000 to 2024-03-11 20:01:41.810000"

# 2. Insert a new row before "Idling time percentage" (originally row 34,
#    now row 35 after the insert above) for "Cycle Count of battery"
$ws.Rows("35:35").Insert()
$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = 136
